{"js": "const lines = [\n  \"\",\n  \"Update: 2026-02-21\",\n  \"- Automation flow builder upgraded with WhatsApp-style canvas node rendering.\",\n  \"- Added richer node palette with icons and expanded WhatsApp bot node options.\",\n  \"- Added Form node support (multi-field config and preview).\",\n  \"- Added template variable mapping validation in builder and publish guard.\",\n  \"- Removed hard quick-reply cap in builder UI as requested.\",\n  \"- Fixed JSX template-label syntax compile failure.\",\n  \"- Fixed CONTENTS panel usability: wider panel, proper scroll, click-safe interactions, no overlap.\"\n];\n\nconst ANCHOR_TEXT = \"4. Full RBAC UI hide/disable sweep for all pages/actions\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that ends the existing \"Completed\" list (the diff\n// inserts the new status block immediately after it, right before the\n// section end).\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text === ANCHOR_TEXT) {\n    anchor = p;\n  }\n}\nif (!anchor) {\n  // Fall back to the last paragraph in the body if the text could not be\n  // matched exactly (e.g. trailing whitespace differences).\n  anchor = paragraphs.items[paragraphs.items.length - 1];\n}\n\nfor (const line of lines) {\n  // insertParagraph on the anchor inherits the \"Helvetica Light\" / 24\n  // half-point run formatting already present on that paragraph.\n  const newPara = anchor.insertParagraph(line, Word.InsertLocation.after);\n  if (line === \"\") {\n    // Force an explicit (empty) run/text node instead of a run with no <w:t>.\n    newPara.insertText(\"\", Word.InsertLocation.start);\n  }\n  anchor = newPara;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$lines = @(\n  \"\",\n  \"Update: 2026-02-21\",\n  \"- Automation flow builder upgraded with WhatsApp-style canvas node rendering.\",\n  \"- Added richer node palette with icons and expanded WhatsApp bot node options.\",\n  \"- Added Form node support (multi-field config and preview).\",\n  \"- Added template variable mapping validation in builder and publish guard.\",\n  \"- Removed hard quick-reply cap in builder UI as requested.\",\n  \"- Fixed JSX template-label syntax compile failure.\",\n  \"- Fixed CONTENTS panel usability: wider panel, proper scroll, click-safe interactions, no overlap.\"\n)\n\n$anchorText = \"4. Full RBAC UI hide/disable sweep for all pages/actions\"\n\n# Locate the paragraph ending the existing \"Completed\" list (the diff inserts\n# the new status block immediately after it, right before the section end).\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq $anchorText) {\n    $anchorIndex = $i\n  }\n}\nif ($anchorIndex -eq -1) {\n  $anchorIndex = $d.Paragraphs.Count\n}\n\n$anchorRange = $d.Paragraphs.Item($anchorIndex).Range\n\nforeach ($line in $lines) {\n  # InsertParagraphAfter() on a range that already carries the\n  # \"Helvetica Light\" / 12pt run formatting propagates that same\n  # formatting to the new paragraph's run automatically.\n  $anchorRange.InsertParagraphAfter()\n  $newPara = $d.Paragraphs.Last\n  $newPara.Range.Text = $line\n  $anchorRange = $newPara.Range\n}\n\nWrite-Host \"Inserted $($lines.Count) paragraphs\"\n"}
